$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.556565037682828
$ws.Range("C2").Value = 0.612606004275778
$ws.Range("L2").Value = 0.599020579536134

$ws.Range("B3").Value = 0.503247309719211
$ws.Range("L3").Value = 0.605713809935328
